$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "274.64"
Set-TextValue $ws.Range("E2") "-1.32%"
Set-TextValue $ws.Range("E3") "-2.03%"
Set-TextValue $ws.Range("D4") "4.901"
Set-TextValue $ws.Range("E4") "1.88%"
Set-TextValue $ws.Range("D5") "0.06325"
Set-TextValue $ws.Range("E5") "1.29%"
Set-TextValue $ws.Range("D6") "6.855"
Set-TextValue $ws.Range("E6") "-0.99%"
Set-TextValue $ws.Range("D7") "3.320"
Set-TextValue $ws.Range("E7") "1.57%"
Set-TextValue $ws.Range("D8") "1.247"
Set-TextValue $ws.Range("E8") "32.45%"
Set-TextValue $ws.Range("D9") "0.8691"
Set-TextValue $ws.Range("E9") "-1.19%"
Set-TextValue $ws.Range("E10") "6.04%"
Set-TextValue $ws.Range("D11") "0.05024"
Set-TextValue $ws.Range("E11") "-4.24%"
Set-TextValue $ws.Range("D12") "0.07476"
Set-TextValue $ws.Range("E12") "1.97%"
Set-TextValue $ws.Range("D13") "0.02954"
Set-TextValue $ws.Range("E13") "-4.23%"
Set-TextValue $ws.Range("E14") "-0.43%"
Set-TextValue $ws.Range("D15") "0.001576"
Set-TextValue $ws.Range("E15") "1.44%"
Set-TextValue $ws.Range("D16") "0.0006319"
Set-TextValue $ws.Range("E16") "0.91%"
Set-TextValue $ws.Range("D17") "0.005851"
Set-TextValue $ws.Range("E17") "-1.01%"
Set-TextValue $ws.Range("D18") "3.445"
Set-TextValue $ws.Range("E18") "-0.20%"
Set-TextValue $ws.Range("D19") "2.271"
Set-TextValue $ws.Range("E19") "-0.57%"
Set-TextValue $ws.Range("E20") "-0.10%"
Set-TextValue $ws.Range("D21") "0.1335"
Set-TextValue $ws.Range("E21") "3.17%"
Set-TextValue $ws.Range("E22") "1.51%"
Set-TextValue $ws.Range("D23") "0.04360"
Set-TextValue $ws.Range("E23") "1.12%"
Set-TextValue $ws.Range("D24") "0.001179"
Set-TextValue $ws.Range("E24") "0.03%"
Set-TextValue $ws.Range("E25") "-0.91%"
Set-TextValue $ws.Range("D26") "0.0001200"
Set-TextValue $ws.Range("E26") "0.08%"
Set-TextValue $ws.Range("D27") "0.0001687"
Set-TextValue $ws.Range("E27") "-0.13%"
Set-TextValue $ws.Range("D40") "0.04043"
Set-TextValue $ws.Range("E40") "0.18%"
Set-TextValue $ws.Range("D41") "0.006690"
Set-TextValue $ws.Range("E41") "-0.21%"
Set-TextValue $ws.Range("D42") "0.1164"
Set-TextValue $ws.Range("E42") "0.88%"
Set-TextValue $ws.Range("D43") "0.002201"
Set-TextValue $ws.Range("E43") "3.20%"
Set-TextValue $ws.Range("D44") "0.01069"
Set-TextValue $ws.Range("E44") "-12.26%"
Set-TextValue $ws.Range("D45") "0.00005308"
Set-TextValue $ws.Range("E45") "4.29%"
Set-TextValue $ws.Range("D47") "1.486"
Set-TextValue $ws.Range("E47") "-37.47%"

Write-Host "Applied symbol list update"
